$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Rows 4-7 correspond to: 36f00e0d..., 45851687..., 4b88e143..., 5e526824...
$rows = 4,5,6,7

foreach ($r in $rows) {
    # Priority column (E): "low" -> "ht" on both locale sheets
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Latest Handoff Datetime column (H) for zh-cn: a fresh handoff was generated
    $zhcn.Range("H$r").Value = "2016-09-05 02:35:28"
}

# The overall "Ready for handoff" timestamp is shared between the Overview sheet
# (column G) and the de-de sheet (column H, Latest Handoff Datetime) for these
# same four rows.
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-05 02:35:33"
    $dede.Range("H$r").Value = "2016-09-05 02:35:33"
}
